# Add a new "Results" slide (slide 8) after the existing slide 7.
#
# Slide 7 is titled "Results" with an empty content placeholder, and is
# reused as the structural basis for the new slide (so the new slide
# inherits the same Title-and-Content layout, group shape transform,
# creationId extLst and colour-map override that a hand-authored slide
# in this deck carries). We then overwrite its title/body text with the
# new "Results" bullet content.

$p = $ppt.ActivePresentation

$s7 = $p.Slides.Item(7)
$s7.Duplicate() | Out-Null

$s8 = $p.Slides.Item(8)

# Title
$s8.Shapes.Item(1).TextFrame.TextRange.Text = "Results"

# Body bullets
$content = $s8.Shapes.Item(2).TextFrame.TextRange
$content.Text = "Significant speedup for large matrices`rStrassen`rParallel algorithm is faster than serial one.`rGranularity changes for different matrix sizes to avoid memory overflow."

# Split "Strassen" into its own run followed by the rest of that bullet,
# matching "Strassen algorithm slower for smaller matrices".
$strassenPara = $content.Paragraphs(2, 1)
$strassenPara.InsertAfter(" algorithm slower for smaller matrices") | Out-Null

# Trailing empty paragraph after the last bullet.
$content.InsertAfter("`r") | Out-Null
